$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "SamplesTab" query (B3) with the new version of the query that
# no longer selects smp.sample_tumor_status ("Tumor") or smp.sample_type
# ("Analyte Type").
$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001437' AND sp.gender = 'Unknown'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value = $newSamplesQuery

# Remove the TsvExcel/WebExcel filename cells from the SamplesTab (row 3) and
# FilesTab (row 4) rows; only the ParticipantsTab row (row 2) keeps them.
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Update the active selection to B3.
$ws.Range("B3").Select() | Out-Null
